# TC03_C3DC_phs002599_Race-Unknown.xlsx
# Fixed query issue for C3DC phs002599
#
# The "TreatmentTab" query (row 5, column B) wrapped its REPLACE() call in a
# redundant CONCAT(...) - strip that so the formula text just reads
# REPLACE(trt.treatment_agent, ';', ', ').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedQuery = @"
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs002599' AND prt.race = 'Unknown'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
"@

# Write the corrected SQL text into the TreatmentTab query cell.
$ws.Range("B5").Value = $fixedQuery

# Re-touch the font on that cell - this mirrors what happened in the
# authored workbook (a fresh style record for B5 distinct from the one
# shared by the other query cells B2/B3/B4/B6/B7).
$ws.Range("B5").Font.ThemeColor = 1
$ws.Range("B5").WrapText = $true

# The author's last save left the grid scrolled/selected near the fixed
# cell rather than at the top of the sheet.
$null = $ws.Range("C5").Select()
